# Daily refresh of the cryptos list (prices + 1h volume deltas).
# Some "Price" values look numeric (e.g. "0.998", "13.57") but must stay
# plain text, matching the source data's inlineStr cells - a leading
# apostrophe forces Excel to keep them as text instead of auto-converting
# to a number, same as typing them in by hand.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.906.42"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").Value = "2.916.51"
$ws.Range("E3").Value = "  -1.78%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'587.66"
$ws.Range("E5").Value = "  -1.60%  "
$ws.Range("D6").Value = "'146.18"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("D9").Value = "2.916.73"
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("D10").Value = "'7.06"
$ws.Range("E10").Value = "  -4.29%  "
$ws.Range("D11").Value = "'0.152"
$ws.Range("E11").Value = "  +5.15%  "
$ws.Range("D12").Value = "'0.438"
$ws.Range("E12").Value = "  -1.94%  "
$ws.Range("D13").Value = "'0.0000240"
$ws.Range("E13").Value = "  +3.51%  "
$ws.Range("E14").Value = "  -3.14%  "
$ws.Range("E15").Value = "  -1.39%  "
$ws.Range("D16").Value = "3.400.66"
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("D17").Value = "61.911.22"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("E18").Value = "  -1.94%  "
$ws.Range("D19").Value = "2.912.90"
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("D20").Value = "'435.69"
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("D21").Value = "'13.57"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("E22").Value = "  -2.52%  "
$ws.Range("D23").Value = "'6.94"
$ws.Range("E23").Value = "  -2.61%  "
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'11.93"
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").Value = "'10.33"
$ws.Range("E26").Value = "  -6.46%  "
$ws.Range("D27").Value = "'2.08"
$ws.Range("E27").Value = "  -4.72%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").Value = "'0.0000110"
$ws.Range("E29").Value = "  +23.62%  "
$ws.Range("E30").Value = "  +4.46%  "
$ws.Range("D31").Value = "'2.57"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("D32").Value = "'2.10"
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("E33").Value = "  +2.30%  "
$ws.Range("D34").Value = "'26.03"
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "'0.976"
$ws.Range("E36").Value = "  -2.51%  "
$ws.Range("E37").Value = "  +5.94%  "
$ws.Range("D38").Value = "'5.53"
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("D39").Value = "'49.13"
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("E41").Value = "  -2.92%  "
$ws.Range("E42").Value = "  -4.25%  "
$ws.Range("E43").Value = "  -2.78%  "
$ws.Range("D44").Value = "'39.27"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").Value = "2.701.71"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").Value = "'134.72"
$ws.Range("E46").Value = "  +0.89%  "
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("D48").Value = "'347.18"
$ws.Range("E48").Value = "  -5.53%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("D51").Value = "'22.41"
$ws.Range("E51").Value = "  -4.14%  "
